# EASY-1978: Make DCT_RIGHTSHOLDER column mandatory
#
# The DCT_RIGHTSHOLDER column (Q) is turned into a "mandatory" column in the
# instructions sheet: its header gets the same yellow highlight used by the
# other mandatory columns, and the example values in the data rows are
# cleared out (while keeping the same "mandatory-yellow" look) so the sheet
# no longer ships with sample values "jan" / "eko" / "linda".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: give it the yellow "mandatory" background, same as the other
# required columns (P1, AF1, AG1, AT1, ...).
$ws.Range("Q1").Interior.Color = 65535

# Data cells: clear the example values and apply the yellow "mandatory"
# background used by other required columns' blank cells.
$ws.Range("Q2").Interior.Color = 65535
$ws.Range("Q2").ClearContents()

$ws.Range("Q3").Interior.Color = 65535
$ws.Range("Q3").ClearContents()

$ws.Range("Q4").Interior.Color = 65535
$ws.Range("Q4").ClearContents()

# Restore the active cell to Q1 (what the header column highlighted).
$ws.Range("Q1").Select()
